# The forest-plot data sheet's "Study" column header is relabeled to
# lowercase "study" (part of reworking the plot to show standardized
# mean differences between treatment arms).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "study"
